$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.193.97"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.834.30"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.82"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6655"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07418"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2936"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.98"
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07758"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").Value = "1.844.55"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.993"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6686"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.97"
$ws.Range("E15").Value = "  -4.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.111"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "29.205.49"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "2.094.98"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.03"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.155"
$ws.Range("E23").Value = "  -2.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9997"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.32"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("E26").Value = "  -2.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.619"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.99"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.514"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.111"
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.044"
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.192"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.865"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7464"
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.649"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").Value = "1.286.23"
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.735"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9352"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.887"
$ws.Range("E42").Value = "  -3.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.08376"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.38"
$ws.Range("E45").Value = "  -2.54%  "
$ws.Range("D46").Value = "1.990.42"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.760"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.06"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("E51").Value = "  -0.93%  "
